# daily auto push: 2026-01-10 18:43 UTC
# Insert two new rows of data before row 609 (continuing the 2026/01/10 土
# block and starting a new 2026/01/11 日 block), shifting all subsequent
# rows down by two. Excel's Rows().Insert() takes care of shifting the
# existing data and extending the used range; we then populate the two
# freshly inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 609:610 - everything from 609 downward shifts
# down by two rows (609->611 ... 650->652).
$ws.Rows("609:610").Insert()

# The date column stores plain text like "2026/12/29", not real Excel
# dates, so force Text format first - otherwise assigning a date-shaped
# string auto-converts it into a date serial number.
$ws.Range("A609:A610").NumberFormat = "@"

# Populate the newly inserted row 609.
$ws.Cells.Item(609, 1).Value = "2026/01/10"
$ws.Cells.Item(609, 2).Value = "土"
$ws.Cells.Item(609, 3).Value = 23
$ws.Cells.Item(609, 4).Value = 201

# Populate the newly inserted row 610.
$ws.Cells.Item(610, 1).Value = "2026/01/11"
$ws.Cells.Item(610, 2).Value = "日"
$ws.Cells.Item(610, 3).Value = 2
$ws.Cells.Item(610, 4).Value = 201
